# Applies the TC_54 edit: rename sheet, tweak a couple of labels, widen the
# "0.000" number format to "###0.000", and insert 8 new historical data rows
# above the existing one (old row 13 becomes row 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab.
$ws.Name = "Data"

# 2. Update a couple of text labels.
$ws.Range("C1").Value = "Real Estate Investment: Residential: Hebei [CONVERTCUR(US Dollars; as reported)]"
$ws.Range("A11").Value = "Function Information"

# 3. Insert 8 blank rows above the existing data row (13), pushing it to 21.
$ws.Rows("13:20").Insert()

# Newly inserted rows pick up the (bold) formatting of row 12 above them;
# the historical data rows use the plain, non-bold default font, so reset it
# before (re)writing values/number formats into rows 13-20.
$ws.Range("A13:C20").Font.Bold = $false

# 4. Populate the newly inserted rows (13-20) plus re-apply values to the
#    row that got pushed down to 21, using the same layout/styles as before.
$data = @(
    @(35765, 3271.34, 393.222199739963),
    @(36495, 6384.48, 771.233913701741),
    @(36861, 7111.74, 859.070071772984),
    @(37226, 8354.02, 1009.30428768043),
    @(37591, 10540.67, 1273.48149401055),
    @(37956, 16421.36, 1983.93092658332),
    @(38322, 22352.63, 2700.63671950512),
    @(38687, 29205.22, 3564.22016109348),
    @(39052, 37962.97, 4762.01847291155)
)

$row = 13
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]

    $ws.Cells.Item($row, 1).NumberFormat = "yyyy"
    $ws.Cells.Item($row, 2).NumberFormat = "###0.000"
    $ws.Cells.Item($row, 3).NumberFormat = "###0.000"

    $row = $row + 1
}
